$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.054.68'
$ws.Range('E2').Value = '  -0.52%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.800.03'
$ws.Range('E3').Value = '  -0.01%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.62'
$ws.Range('E5').Value = '  -1.59%  '
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5086'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3856'
$ws.Range('E8').Value = '  +1.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07695'
$ws.Range('E9').Value = '  -3.76%  '
$ws.Range('E10').Value = '  -0.21%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '40.73'
$ws.Range('E11').Value = '  -1.57%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.331'
$ws.Range('E12').Value = '  +0.15%  '
$ws.Range('E13').Value = '  -0.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.30'
$ws.Range('E14').Value = '  -1.57%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.799.79'
$ws.Range('E15').Value = '  -0.25%  '
$ws.Range('E16').Value = '  -0.77%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '92.09'
$ws.Range('E17').Value = '  -0.94%  '
$ws.Range('E18').Value = '  -2.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06571'
$ws.Range('E19').Value = '  -0.65%  '
$ws.Range('E20').Value = '  -0.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.23'
$ws.Range('E21').Value = '  -0.62%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.964'
$ws.Range('E22').Value = '  -0.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.082.23'
$ws.Range('E23').Value = '  -0.58%  '
$ws.Range('E24').Value = '  -1.36%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.240'
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '159.85'
$ws.Range('E26').Value = '  +0.29%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.009.97'
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.406'
$ws.Range('E28').Value = '  +0.76%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.23'
$ws.Range('E29').Value = '  -1.41%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '127.20'
$ws.Range('E30').Value = '  +3.45%  '
$ws.Range('E31').Value = '  -0.92%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.044'
$ws.Range('E32').Value = '  -1.63%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.649'
$ws.Range('E33').Value = '  -0.33%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.533'
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.06967'
$ws.Range('E35').Value = '  -4.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.057'
$ws.Range('E36').Value = '  +2.54%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02331'
$ws.Range('E37').Value = '  +0.93%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2162'
$ws.Range('E38').Value = '  -0.13%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.995'
$ws.Range('E39').Value = '  -1.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '11.44'
$ws.Range('E40').Value = '  -6.30%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6098'
$ws.Range('E41').Value = '  -1.54%  '
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.148'
$ws.Range('E43').Value = '  -1.54%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.22'
$ws.Range('E44').Value = '  +0.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.299'
$ws.Range('E45').Value = '  -5.12%  '
$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.711'
$ws.Range('E46').Value = '  -1.23%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5881'
$ws.Range('E47').Value = '  -1.65%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '125.42'
$ws.Range('E48').Value = '  -1.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.186'
$ws.Range('E49').Value = '  -1.30%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.920'
$ws.Range('E50').Value = '  -0.43%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06726'
$ws.Range('E51').Value = '  -1.50%  '
